# "big update fix bug"
# Adds a "Note" tracking column that marks test cases as Done/DOne.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("Test cases"): mark the "extra /frontend in links" bug row as Done
$ws1 = $wb.Worksheets.Item("Test cases")
$ws1.Range("G13").Value = "Done"

# --- Sheet 2 ("TODOs"): add a "Note" column (E) tracking completion status
$ws2 = $wb.Worksheets.Item("TODOs")
$ws2.Range("E1").Value = "Note"
$ws2.Range("E3").Value  = "Done"
$ws2.Range("E4").Value  = "Done"
$ws2.Range("E5").Value  = "Done"
$ws2.Range("E6").Value  = "Done"
$ws2.Range("E7").Value  = "Done"
$ws2.Range("E9").Value  = "Done"
$ws2.Range("E10").Value = "Done"
$ws2.Range("E11").Value = "DOne"
$ws2.Range("E12").Value = "Done"
$ws2.Range("E13").Value = "Done"
$ws2.Range("E14").Value = "Done"

# Widen column C on the TODOs sheet to fit the existing text now that layout changed
$ws2.Columns.Item(3).ColumnWidth = 35.1666666666666

# Restore cursor positions as left by the author
[void]$ws1.Range("D27").Select()
[void]$ws2.Activate()
[void]$ws2.Range("G19").Select()
